# Update the valid-login test email on the ValidLoginsheet tab
$wb = $excel.ActiveWorkbook

$wsValidLogin = $wb.Worksheets.Item("ValidLoginsheet")
$wsValidLogin.Range("A2").Value = "adityapawar123@yopmail.com"

# Move to the signup data sheet and leave the selection on C2 (last worked-on cell)
$wsSignup = $wb.Worksheets.Item("signupdata")
$wsSignup.Activate()
$wsSignup.Range("C2").Select()
